$d = $word.ActiveDocument

# Locate the final (empty) paragraph of the document - it must stay untouched.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$r = $lastPara.Range
$r.Collapse(0)

# Create 5 fresh empty paragraphs after it:
#   +1 -> stays empty (new "<w:p/>" from the diff)
#   +2 -> will hold the "-note that two communities..." text
#   +3 -> will hold the "- Visualize just with a  subsample" text
#   +4 -> will hold the "-TopK (topD) can be set to > 1 to have " text
#   +5 -> throwaway placeholder, removed again at the end (keeps the
#         replace-operations below from growing the very last paragraph
#         of the story, which behaves differently from interior ones)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.InsertParagraphAfter()

function Fill-Paragraph($paraIndex, $innerXml) {
    $para = $d.Paragraphs.Item($paraIndex)
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           '<w:p>' + $innerXml + '</w:p>' +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

$base = $lastIndex

$noteXml = '<w:r><w:t xml:space="preserve">-note that two communities (e.g. 2 and 398) can have a very similar topic representation. They just happen </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>tobe</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> constructed with different nodes. Nodes that were active while community 2 was relevant are most likely gone when 398 nodes are getting relevant (sliding window approach)</w:t></w:r>'
Fill-Paragraph ($base + 2) $noteXml

$visualizeXml = '<w:r><w:lastRenderedPageBreak/><w:t>- Visualize just with a  subsample</w:t></w:r>'
Fill-Paragraph ($base + 3) $visualizeXml

$topkXml = '<w:r><w:t>-</w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>TopK</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/><w:r><w:t>topD</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
           '<w:r><w:t xml:space="preserve">) can be set to &gt; 1 to have </w:t></w:r>'
Fill-Paragraph ($base + 4) $topkXml

# Drop the throwaway placeholder (+5) that is still empty and unused.
$extraIndex = $d.Paragraphs.Count
$extraPara = $d.Paragraphs.Item($extraIndex)
$extraPara.Range.Delete()
